$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = [double]"1"
$ws.Range("F2").Value = [double]"0.3333333333333333"
$ws.Range("G2").Value = [double]"0.01675466666666667"
$ws.Range("H2").Value = [double]"0.050264"
$ws.Range("I2").Value = [double]"0.0001854906931657378"
$ws.Range("J2").Value = [double]"0.0001854906931657378"
$ws.Range("K2").Value = [double]"2"
$ws.Range("L2").Value = [double]"0.6666666666666666"
$ws.Range("M2").Value = [double]"0.1465046666666667"
$ws.Range("N2").Value = [double]"0.439514"
$ws.Range("O2").Value = [double]"0.07745172725947863"
$ws.Range("P2").Value = [double]"0.07745172725947864"
$ws.Range("Q2").Value = [double]"0.002454636855111111"
$ws.Range("R2").Value = [double]"0.022091731696"
$ws.Range("S2").Value = [double]"1.436657457624436E-05"
$ws.Range("T2").Value = [double]"1.436657457624436E-05"
$ws.Range("E3").Value = [double]"1"
$ws.Range("F3").Value = [double]"0.3333333333333333"
$ws.Range("G3").Value = [double]"0.01675466666666667"
$ws.Range("H3").Value = [double]"0.050264"
$ws.Range("I3").Value = [double]"0.0001854906931657378"
$ws.Range("J3").Value = [double]"0.0001854906931657378"
$ws.Range("N3").Value = [double]"4.707498"
$ws.Range("O3").Value = [double]"0.8295614045753745"
$ws.Range("P3").Value = [double]"0.8295614045753745"
$ws.Range("Q3").Value = [double]"0.02629085327466667"
$ws.Range("R3").Value = [double]"0.236617679472"
$ws.Range("S3").Value = [double]"0.0001538759199582293"
$ws.Range("T3").Value = [double]"0.0001538759199582292"
$ws.Range("E4").Value = [double]"1"
$ws.Range("F4").Value = [double]"0.3333333333333333"
$ws.Range("G4").Value = [double]"0.01675466666666667"
$ws.Range("H4").Value = [double]"0.050264"
$ws.Range("I4").Value = [double]"0.0001854906931657378"
$ws.Range("J4").Value = [double]"0.0001854906931657378"
$ws.Range("O4").Value = [double]"0.09298686816514685"
$ws.Range("P4").Value = [double]"0.09298686816514684"
$ws.Range("Q4").Value = [double]"0.002946983904888889"
$ws.Range("R4").Value = [double]"0.026522855144"
$ws.Range("S4").Value = [double]"1.724819863126417E-05"
$ws.Range("T4").Value = [double]"1.724819863126416E-05"
$ws.Range("I5").Value = [double]"0.9933938536206305"
$ws.Range("J5").Value = [double]"0.9933938536206304"
$ws.Range("K5").Value = [double]"2"
$ws.Range("L5").Value = [double]"0.6666666666666666"
$ws.Range("M5").Value = [double]"0.1465046666666667"
$ws.Range("N5").Value = [double]"0.439514"
$ws.Range("O5").Value = [double]"0.07745172725947863"
$ws.Range("P5").Value = [double]"0.07745172725947864"
$ws.Range("Q5").Value = [double]"13.14578711805933"
$ws.Range("R5").Value = [double]"118.312084062534"
$ws.Range("S5").Value = [double]"0.07694006981186752"
$ws.Range("T5").Value = [double]"0.07694006981186752"
$ws.Range("I6").Value = [double]"0.9933938536206305"
$ws.Range("J6").Value = [double]"0.9933938536206304"
$ws.Range("N6").Value = [double]"4.707498"
$ws.Range("O6").Value = [double]"0.8295614045753745"
$ws.Range("P6").Value = [double]"0.8295614045753745"
$ws.Range("S6").Value = [double]"0.8240812005060743"
$ws.Range("T6").Value = [double]"0.8240812005060741"
$ws.Range("I7").Value = [double]"0.9933938536206305"
$ws.Range("J7").Value = [double]"0.9933938536206304"
$ws.Range("O7").Value = [double]"0.09298686816514685"
$ws.Range("P7").Value = [double]"0.09298686816514684"
$ws.Range("S7").Value = [double]"0.09237258330268876"
$ws.Range("T7").Value = [double]"0.09237258330268873"
$ws.Range("G8").Value = [double]"0.5799533333333334"
$ws.Range("I8").Value = [double]"0.006420655686203657"
$ws.Range("J8").Value = [double]"0.006420655686203655"
$ws.Range("K8").Value = [double]"2"
$ws.Range("L8").Value = [double]"0.6666666666666666"
$ws.Range("M8").Value = [double]"0.1465046666666667"
$ws.Range("N8").Value = [double]"0.439514"
$ws.Range("O8").Value = [double]"0.07745172725947863"
$ws.Range("P8").Value = [double]"0.07745172725947864"
$ws.Range("Q8").Value = [double]"0.08496586978222224"
$ws.Range("R8").Value = [double]"0.7646928280400002"
$ws.Range("S8").Value = [double]"0.0004972908730348662"
$ws.Range("T8").Value = [double]"0.0004972908730348662"
$ws.Range("G9").Value = [double]"0.5799533333333334"
$ws.Range("I9").Value = [double]"0.006420655686203657"
$ws.Range("J9").Value = [double]"0.006420655686203655"
$ws.Range("N9").Value = [double]"4.707498"
$ws.Range("O9").Value = [double]"0.8295614045753745"
$ws.Range("P9").Value = [double]"0.8295614045753745"
$ws.Range("Q9").Value = [double]"0.9100430522533335"
$ws.Range("R9").Value = [double]"8.190387470280001"
$ws.Range("S9").Value = [double]"0.005326328149341971"
$ws.Range("T9").Value = [double]"0.005326328149341969"
$ws.Range("G10").Value = [double]"0.5799533333333334"
$ws.Range("I10").Value = [double]"0.006420655686203657"
$ws.Range("J10").Value = [double]"0.006420655686203655"
$ws.Range("O10").Value = [double]"0.09298686816514685"
$ws.Range("P10").Value = [double]"0.09298686816514684"
$ws.Range("R10").Value = [double]"0.9180736660600001"
$ws.Range("S10").Value = [double]"0.00059703666382682"
$ws.Range("T10").Value = [double]"0.0005970366638268197"
